# Apply updated market-price figures (currentAveragePrice / Leve profit columns)
# pulled by the scheduled pricing runner, across the affected Leve rows on each job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 2995.5854
$ws.Range("I15").Value = 2995.5854
$ws.Range("K15").Value = 8986.7562
$ws.Range("M15").Value = -8817.7562
# Row 98
$ws.Range("H98").Value = 1712.2174
$ws.Range("I98").Value = 1303.8572
$ws.Range("K98").Value = 1303.8572
$ws.Range("M98").Value = 194.1428000000001
# Row 106
$ws.Range("H106").Value = 1250
# Row 107
$ws.Range("H107").Value = 809.0741
$ws.Range("I107").Value = 832.13043
$ws.Range("J107").Value = 676.5
$ws.Range("K107").Value = 832.13043
$ws.Range("L107").Value = 676.5
$ws.Range("M107").Value = 1087.86957
$ws.Range("N107").Value = -4516.5
# Row 111
$ws.Range("H111").Value = 951.7143
$ws.Range("I111").Value = 907.25
$ws.Range("J111").Value = 1011
$ws.Range("K111").Value = 2721.75
$ws.Range("L111").Value = 3033
$ws.Range("M111").Value = 345.25
$ws.Range("N111").Value = -9167
# Row 122
$ws.Range("H122").Value = 1712.2174
$ws.Range("I122").Value = 1303.8572
$ws.Range("K122").Value = 3911.5716
$ws.Range("M122").Value = -1461.5716

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7195.5444
$ws.Range("I32").Value = 5884.311
$ws.Range("K32").Value = 5884.311
$ws.Range("M32").Value = -5597.311
# Row 61
$ws.Range("H61").Value = 3032.35
$ws.Range("I61").Value = 487.15384
$ws.Range("J61").Value = 7759.143
$ws.Range("K61").Value = 487.15384
$ws.Range("L61").Value = 7759.143
$ws.Range("M61").Value = -275.15384
$ws.Range("N61").Value = -8183.143
# Row 136
$ws.Range("H136").Value = 3032.35
$ws.Range("I136").Value = 487.15384
$ws.Range("J136").Value = 7759.143
$ws.Range("K136").Value = 1461.46152
$ws.Range("L136").Value = 23277.429
$ws.Range("M136").Value = 1088.53848
$ws.Range("N136").Value = -28377.429

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 2955.138
$ws.Range("I134").Value = 2028.5834
$ws.Range("J134").Value = 7402.6
$ws.Range("K134").Value = 6085.7502
$ws.Range("L134").Value = 22207.8
$ws.Range("M134").Value = -3550.7502
$ws.Range("N134").Value = -27277.8

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2124.9167
$ws.Range("I16").Value = 1340
$ws.Range("J16").Value = 2685.5715
$ws.Range("K16").Value = 1340
$ws.Range("L16").Value = 2685.5715
$ws.Range("M16").Value = -1053
$ws.Range("N16").Value = -3259.5715
# Row 31
$ws.Range("H31").Value = 2085576.6
$ws.Range("I31").Value = 2501633.8
$ws.Range("J31").Value = 5291
$ws.Range("K31").Value = 2501633.8
$ws.Range("L31").Value = 5291
$ws.Range("M31").Value = -2501338.8
$ws.Range("N31").Value = -5881
# Row 34
$ws.Range("H34").Value = 2085576.6
$ws.Range("I34").Value = 2501633.8
$ws.Range("J34").Value = 5291
$ws.Range("K34").Value = 2501633.8
$ws.Range("L34").Value = 5291
$ws.Range("M34").Value = -2501431.8
$ws.Range("N34").Value = -5695
# Row 113
$ws.Range("H113").Value = 2124.9167
$ws.Range("I113").Value = 1340
$ws.Range("J113").Value = 2685.5715
$ws.Range("K113").Value = 1340
$ws.Range("L113").Value = 2685.5715
$ws.Range("M113").Value = 830
$ws.Range("N113").Value = -7025.5715

$ws = $wb.Worksheets.Item("CUL")
# Row 97
$ws.Range("H97").Value = 2083.5454
$ws.Range("I97").Value = 799.5
$ws.Range("J97").Value = 2368.889
$ws.Range("K97").Value = 2398.5
$ws.Range("L97").Value = 7106.667
$ws.Range("M97").Value = -1902.5
$ws.Range("N97").Value = -8098.667
# Row 98
$ws.Range("H98").Value = 200
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0

$ws = $wb.Worksheets.Item("GSM")
# Row 45
$ws.Range("H45").Value = 30000
$ws.Range("J45").Value = 30000
$ws.Range("L45").Value = 30000
$ws.Range("N45").Value = -31118
# Row 102
$ws.Range("H102").Value = 30090.139
$ws.Range("I102").Value = 1731.5
$ws.Range("J102").Value = 86807.414
$ws.Range("K102").Value = 1731.5
$ws.Range("L102").Value = 86807.414
$ws.Range("M102").Value = -109.5
$ws.Range("N102").Value = -90051.414
# Row 122
$ws.Range("H122").Value = 3737.0417
$ws.Range("I122").Value = 2321.1538
$ws.Range("J122").Value = 5410.364
$ws.Range("K122").Value = 6963.4614
$ws.Range("L122").Value = 16231.092
$ws.Range("M122").Value = -4513.4614
$ws.Range("N122").Value = -21131.092
# Row 132
$ws.Range("H132").Value = 3353.561
$ws.Range("I132").Value = 3308.25
$ws.Range("J132").Value = 3417.5293
$ws.Range("K132").Value = 9924.75
$ws.Range("L132").Value = 10252.5879
$ws.Range("M132").Value = -7394.75
$ws.Range("N132").Value = -15312.5879

$ws = $wb.Worksheets.Item("LTW")
# Row 3
$ws.Range("H3").Value = 61670.832
$ws.Range("J3").Value = 61670.832
$ws.Range("L3").Value = 61670.832
$ws.Range("N3").Value = -61894.832
# Row 7
$ws.Range("H7").Value = 2146.2666
$ws.Range("I7").Value = 1213.5
$ws.Range("J7").Value = 2768.111
$ws.Range("K7").Value = 1213.5
$ws.Range("L7").Value = 2768.111
$ws.Range("M7").Value = -1101.5
$ws.Range("N7").Value = -2992.111
# Row 15
$ws.Range("H15").Value = 61670.832
$ws.Range("J15").Value = 61670.832
$ws.Range("L15").Value = 61670.832
$ws.Range("N15").Value = -62010.832
# Row 55
$ws.Range("H55").Value = 1224.5385
$ws.Range("I55").Value = 180.83333
$ws.Range("J55").Value = 2119.1428
$ws.Range("K55").Value = 180.83333
$ws.Range("L55").Value = 2119.1428
$ws.Range("M55").Value = -7.833329999999989
$ws.Range("N55").Value = -2465.1428
# Row 100
$ws.Range("H100").Value = 1944.0555
$ws.Range("I100").Value = 1315.3
$ws.Range("J100").Value = 2730
$ws.Range("K100").Value = 1315.3
$ws.Range("L100").Value = 2730
$ws.Range("M100").Value = -774.3
$ws.Range("N100").Value = -3812
# Row 126
$ws.Range("H126").Value = 2146.2666
$ws.Range("I126").Value = 1213.5
$ws.Range("J126").Value = 2768.111
$ws.Range("K126").Value = 3640.5
$ws.Range("L126").Value = 8304.332999999999
$ws.Range("M126").Value = -1170.5
$ws.Range("N126").Value = -13244.333

$ws = $wb.Worksheets.Item("WVR")
# Row 82
$ws.Range("H82").Value = 30622.846
$ws.Range("J82").Value = 30622.846
$ws.Range("L82").Value = 30622.846
$ws.Range("N82").Value = -31388.846
# Row 85
$ws.Range("H85").Value = 30622.846
$ws.Range("J85").Value = 30622.846
$ws.Range("L85").Value = 30622.846
$ws.Range("N85").Value = -33274.84600000001

# CUL row 98: N98 (LeveProfitHQ) no longer applicable for this update - remove it
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N98").ClearContents()

Write-Output "Applied scheduled price updates."
